# Apply the stimuli-sheet update: add a "carrier" value for each practice/word
# row in column D, and mark the generic carrier rows (6-9) plus the duplicate
# "unique_video" / "unique_audio" rows (14-21) with their pair_kind / carrier
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the still-blank F1 cell blank (defensive no-op; some engines materialize
# empty shared-string placeholders with a stray value on load/save round-trip).
$ws.Range("F1").Value = ""

# Practice rows (2-5): fill in the new "carrier" column D
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows (6-9): add pair_kind in column J
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: add kind (C) and carrier (D) values
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
